$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = '59.141.73'
$ws.Range("E2").Value = '  +3.29%  '

$ws.Range("D3").Value = '3.119.32'
$ws.Range("E3").Value = '  +1.41%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  -0.04%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '524.88'
$ws.Range("E5").Value = '  +1.88%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '145.47'
$ws.Range("E6").Value = '  +3.15%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("E7").Value = '  -0.04%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.440'
$ws.Range("E8").Value = '  +1.34%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '7.43'
$ws.Range("E9").Value = '  +2.52%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.110'
$ws.Range("E10").Value = '  +1.31%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.385'
$ws.Range("E11").Value = '  +3.60%  '

$ws.Range("D12").Value = '3.649.40'
$ws.Range("E12").Value = '  +1.47%  '

$ws.Range("E13").Value = '  +1.59%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '27.34'
$ws.Range("E14").Value = '  +7.33%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0000168'
$ws.Range("E15").Value = '  +2.35%  '

$ws.Range("D16").Value = '59.065.84'
$ws.Range("E16").Value = '  +3.08%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '6.25'
$ws.Range("E17").Value = '  +5.72%  '

$ws.Range("D18").Value = '3.121.67'
$ws.Range("E18").Value = '  +1.53%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.11'
$ws.Range("E19").Value = '  +0.66%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '8.32'
$ws.Range("E20").Value = '  +1.92%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '341.11'
$ws.Range("E21").Value = '  +1.20%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.999'
$ws.Range("E22").Value = '  -0.26%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.512'
$ws.Range("E23").Value = '  +2.17%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '66.04'

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.173'
$ws.Range("E25").Value = '  +1.91%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.00'
$ws.Range("E26").Value = '  -0.10%  '

$ws.Range("D27").Value = '0.0₃0932'
$ws.Range("E27").Value = '  -2.00%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.69'
$ws.Range("E28").Value = '  +3.19%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.37'
$ws.Range("E29").Value = '  +3.97%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.85'
$ws.Range("E30").Value = '  +2.29%  '

$ws.Range("E31").Value = '  +4.07%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '21.17'
$ws.Range("E32").Value = '  +2.08%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '155.42'
$ws.Range("E33").Value = '  +0.68%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.68'
$ws.Range("E34").Value = '  +3.10%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '6.18'
$ws.Range("E35").Value = '  +5.57%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '27.25'
$ws.Range("E36").Value = '  +3.93%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.31'
$ws.Range("E37").Value = '  +6.05%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0688'
$ws.Range("E38").Value = '  +2.46%  '

$ws.Range("E39").Value = '  +3.22%  '

$ws.Range("D40").Value = '3.156.80'
$ws.Range("E40").Value = '  +1.36%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '36.97'
$ws.Range("E41").Value = '  -0.05%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.999'
$ws.Range("E42").Value = '  -0.02%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.666'
$ws.Range("E43").Value = '  -0.70%  '

$ws.Range("B44").Value = 'Stacks'
$ws.Range("C44").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.46'
$ws.Range("E44").Value = '  +5.37%  '

$ws.Range("B45").Value = 'Maker'
$ws.Range("C45").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D45").Value = '2.292.57'
$ws.Range("E45").Value = '  +2.37%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0259'
$ws.Range("E46").Value = '  +2.34%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '21.13'
$ws.Range("E47").Value = '  +5.64%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.962'
$ws.Range("E48").Value = '  +1.45%  '

$ws.Range("E49").Value = '  +3.00%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.753'
$ws.Range("E50").Value = '  +9.66%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '261.82'
$ws.Range("E51").Value = '  +11.26%  '
